$wb = $excel.ActiveWorkbook

# --- Content change: localization status "Ready for handoff" -> "In Translation" ---
# Overview sheet: zh-cn / de-de status columns (E, F) for both data rows.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F3").Value = "In Translation"

# zh-cn sheet: Status column (C) for both data rows.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C3").Value = "In Translation"

# de-de sheet: Status column (C) for both data rows.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C3").Value = "In Translation"

# --- Column width adjustments that follow from the shorter status text ---
# (these columns got narrower once "Ready for handoff" became "In Translation")
$wsOverview.Columns("E:F").ColumnWidth = 12.5
$wsZhCn.Columns("C").ColumnWidth = 12.5
$wsDeDe.Columns("C").ColumnWidth = 12.5
